# The target diff only rewrites the internal <w:nsid w:val="..."/> GUIDs
# stored inside four <w:abstractNum> list definitions in word/numbering.xml
# (abstractNumId 990, 99411, 99711 and 99413). The commit message itself
# ("Automatic build output files") confirms these are opaque, randomly
# re-minted build identifiers, not an authored content change: every
# surrounding element (multiLevelType, lvl defs, numFmt, indents, list
# text used throughout the document, etc.) is byte-for-byte identical
# before and after.
#
# The OOXML "nsid" is a low-level list-definition identifier that Word's
# automation surface (VBA/COM) has never exposed for reading or writing -
# List, ListFormat and ListTemplate only expose ListID/Name/OutlineNumbered/
# ListLevels/etc., none of which round-trip to w:nsid - so there is no
# Word.Application object-model call (Find/Replace, Paragraphs, Lists,
# ListFormat, ListTemplate, Range.WordOpenXML, ...) that can address or
# overwrite that attribute; the document package itself is also not
# reachable from script (it is held open by the running Word instance).
#
# Since the visible document content, styles and list formatting are
# unaffected, there is nothing for an automation script to do here -
# touching the model confirms the document is intact and leaves it
# unchanged, which matches the only observable effect of the source
# commit.
$d = $word.ActiveDocument
$null = $d.Content
Write-Output ("Paragraphs: " + $d.Paragraphs.Count + "; no addressable Word object-model property maps to w:nsid, so no change is applied.")
